$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its text representation instead of being
# auto-converted to numbers by the COM Value setter (many prices look numeric,
# e.g. "1.00", "586.69"). Temporarily force Text format, write the values, then
# restore the default "Normal" style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.064.14"
$ws.Range("E2").Value = "  +4.61%  "
$ws.Range("D3").Value = "3.462.40"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "585.66"
$ws.Range("E5").Value = "  +6.33%  "
$ws.Range("D6").Value = "186.87"
$ws.Range("E6").Value = "  +8.42%  "
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").Value = "3.456.09"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").Value = "0.645"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "56.32"
$ws.Range("E12").Value = "  +6.17%  "
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "9.39"
$ws.Range("E14").Value = "  +3.99%  "
$ws.Range("D15").Value = "4.023.44"
$ws.Range("E15").Value = "  +4.37%  "
$ws.Range("D16").Value = "18.72"
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("D17").Value = "3.471.93"
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("D18").Value = "67.041.29"
$ws.Range("E18").Value = "  +4.77%  "
$ws.Range("D19").Value = "12.14"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "1.01"
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").Value = "489.84"
$ws.Range("E22").Value = "  +9.82%  "
$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +8.24%  "
$ws.Range("E24").Value = "  +22.60%  "
$ws.Range("D25").Value = "4.45"
$ws.Range("E25").Value = "  +10.63%  "
$ws.Range("D26").Value = "89.70"
$ws.Range("E26").Value = "  +3.54%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("D28").Value = "10.93"
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("D29").Value = "9.05"
$ws.Range("E29").Value = "  +6.06%  "
$ws.Range("D30").Value = "31.34"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  +9.84%  "
$ws.Range("D32").Value = "599.36"
$ws.Range("E32").Value = "  +5.51%  "
$ws.Range("D33").Value = "11.70"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").Value = "64.02"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E35").Value = "  +4.98%  "
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "36.40"
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("D39").Value = "3.53"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "0.383"
$ws.Range("E40").Value = "  +4.86%  "
$ws.Range("D41").Value = "0.0₃0753"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("D42").Value = "3.231.50"
$ws.Range("E42").Value = "  +5.92%  "
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +6.72%  "
$ws.Range("D44").Value = "0.0428"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.78"
$ws.Range("E46").Value = "  +24.64%  "
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +13.17%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "8.69"
$ws.Range("E51").Value = "  +6.86%  "

$priceRange.Style = "Normal"
